$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Siege Analytics bullet: split the run so "50M" is bold/colored
#    like the other inline stat call-outs in this bullet. Do this
#    BEFORE the professional-summary replacement below so the only
#    "50M" in the document at this point is the one to be bolded.
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic race coding errors affecting 50M voters, developed geospatial machine",
    2) | Out-Null

$boldRng = $d.Content
$boldRng.Find.Execute("50M", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$boldRng.Font.Bold = 1
$boldRng.Font.Color = 5258796

# -----------------------------------------------------------------
# 2) Professional summary paragraph: "affecting all Black and
#    Asian-American voters" -> "affecting 50M voters" (plain text).
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed geospatial ML",
    2) | Out-Null

# -----------------------------------------------------------------
# 3) Key Projects impact line: "affecting all Black and
#    Asian-American voters" -> "affecting 50M voters nationwide".
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved electoral prediction accuracy by 22%",
    2) | Out-Null

# -----------------------------------------------------------------
# 4) Move the "Software Engineer - Mautinoa Technologies" job block
#    (heading + 4 paragraphs) so it follows the "Partner - Siege
#    Analytics" block instead of the "Software Engineer - Salsa
#    Labs" block.
# -----------------------------------------------------------------
$startRng = $d.Content
$startRng.Find.Execute("Software Engineer - Mautinoa Technologies", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pStart = $startRng.Paragraphs(1)

$endRng = $d.Content
$endRng.Find.Execute("Geospatial analysis on populations and boundaries for impact assessment", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pEnd = $endRng.Paragraphs(1)

$blockStart = $pStart.Range.Start
$blockEnd = $pEnd.Range.End
$block = $d.Range($blockStart, $blockEnd)
$blockFT = $block.FormattedText

$target = $d.Content
$target.Find.Execute("Senior Analyst - Myers Research", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pTarget = $target.Paragraphs(1)
$insertStart = $pTarget.Range.Start

$insertPoint = $d.Range($insertStart, $insertStart)
$insertPoint.FormattedText = $blockFT

# Re-apply the Heading3 style to the relocated job-title paragraph
# (FormattedText insertion at a collapsed range does not carry the
# paragraph style of the incoming paragraph mark).
$headingFix = $d.Range($insertStart, $insertStart)
$headingFix.Find.Execute("Software Engineer - Mautinoa Technologies", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingPara = $headingFix.Paragraphs(1)
$headingPara.Range.Style = "Heading3"

# Remove the original copy of the block, now shifted later in the
# document by the length of the text we just inserted.
$shift = $blockEnd - $blockStart
$origBlock = $d.Range($blockStart + $shift, $blockEnd + $shift)
$origBlock.Delete() | Out-Null

Write-Output "edit complete"
